# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across the
# Leve market-data tables on each job sheet to reflect refreshed market board
# prices pulled in by the scheduled data-refresh runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1) # ALC
# Row 18: You Grow, Girl / Growth Formula Beta
$ws.Cells.Item(18, 8).Value = 5237.7
$ws.Cells.Item(18, 9).Value = 3999
$ws.Cells.Item(18, 11).Value = 3999
$ws.Cells.Item(18, 13).Value = -3715

# Row 33: Glazed and Confused / Clear Glass Lens
$ws.Cells.Item(33, 8).Value = 304.26315
$ws.Cells.Item(33, 9).Value = 337.75
$ws.Cells.Item(33, 11).Value = 337.75
$ws.Cells.Item(33, 13).Value = -108.75

# Row 43: Growing Is Knowing / Growth Formula Gamma
$ws.Cells.Item(43, 8).Value = 6313
$ws.Cells.Item(43, 9).Value = 8000.5
$ws.Cells.Item(43, 11).Value = 8000.5
$ws.Cells.Item(43, 13).Value = -7931.5

# Row 88: The Grave of Hemlock Groves / Growth Formula Zeta
$ws.Cells.Item(88, 8).Value = 2683.625
$ws.Cells.Item(88, 9).Value = 2890
$ws.Cells.Item(88, 10).Value = 2614.8333
$ws.Cells.Item(88, 11).Value = 2890
$ws.Cells.Item(88, 12).Value = 2614.8333
$ws.Cells.Item(88, 13).Value = -2484
$ws.Cells.Item(88, 14).Value = -3426.8333

# Row 91: Dappling the Highlands (L) / Growth Formula Zeta
$ws.Cells.Item(91, 8).Value = 2683.625
$ws.Cells.Item(91, 9).Value = 2890
$ws.Cells.Item(91, 10).Value = 2614.8333
$ws.Cells.Item(91, 11).Value = 2890
$ws.Cells.Item(91, 12).Value = 2614.8333
$ws.Cells.Item(91, 13).Value = -1486
$ws.Cells.Item(91, 14).Value = -5422.8333

# Row 129: Practical Command / Commanding Craftsman's Draught
$ws.Cells.Item(129, 8).Value = 486.875
$ws.Cells.Item(129, 9).Value = 486.875
$ws.Cells.Item(129, 11).Value = 1460.625
$ws.Cells.Item(129, 13).Value = 3539.375

# Row 135: For Tired Minds / Grade 1 Gemsap of Intelligence
$ws.Cells.Item(135, 8).Value = 1686
$ws.Cells.Item(135, 9).Value = 1686
$ws.Cells.Item(135, 11).Value = 15174
$ws.Cells.Item(135, 13).Value = -12639

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Cells.Item(137, 8).Value = 4821
$ws.Cells.Item(137, 10).Value = 6237
$ws.Cells.Item(137, 12).Value = 18711
$ws.Cells.Item(137, 14).Value = -23811

$ws = $wb.Worksheets.Item(2) # ARM
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Cells.Item(2, 8).Value = 924.1539
$ws.Cells.Item(2, 9).Value = 924.1539
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 924.1539
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = -811.1539
$ws.Cells.Item(2, 14).ClearContents()

# Row 98: Greaving / Doman Iron Greaves of Maiming
$ws.Cells.Item(98, 8).Value = 15000
$ws.Cells.Item(98, 10).Value = 15000
$ws.Cells.Item(98, 12).Value = 15000
$ws.Cells.Item(98, 14).Value = -20990

# Row 116: No Scope / Titanbronze Ingot
$ws.Cells.Item(116, 8).Value = 924.1539
$ws.Cells.Item(116, 9).Value = 924.1539
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 924.1539
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 13).Value = 1369.8461
$ws.Cells.Item(116, 14).ClearContents()

$ws = $wb.Worksheets.Item(3) # BSM
# Row 3: Hells Bells / Bronze Ingot
$ws.Cells.Item(3, 8).Value = 924.1539
$ws.Cells.Item(3, 9).Value = 924.1539
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 924.1539
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = -810.1539
$ws.Cells.Item(3, 14).ClearContents()

# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Cells.Item(86, 8).Value = 2514.1904
$ws.Cells.Item(86, 9).Value = 2755.5
$ws.Cells.Item(86, 10).Value = 1066.3334
$ws.Cells.Item(86, 11).Value = 2755.5
$ws.Cells.Item(86, 12).Value = 1066.3334
$ws.Cells.Item(86, 13).Value = -1632.5
$ws.Cells.Item(86, 14).Value = -3312.3334

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Cells.Item(89, 8).Value = 2514.1904
$ws.Cells.Item(89, 9).Value = 2755.5
$ws.Cells.Item(89, 10).Value = 1066.3334
$ws.Cells.Item(89, 11).Value = 13777.5
$ws.Cells.Item(89, 12).Value = 5331.666999999999
$ws.Cells.Item(89, 13).Value = -8161.5
$ws.Cells.Item(89, 14).Value = -16563.667

$ws = $wb.Worksheets.Item(4) # CRP
# Row 86: Birch, Please / Birch Lumber
$ws.Cells.Item(86, 8).Value = 6271.143
$ws.Cells.Item(86, 9).Value = 4974.75
$ws.Cells.Item(86, 11).Value = 4974.75
$ws.Cells.Item(86, 13).Value = -3851.75

# Row 89: Built This City on Blocks and Soul (L) / Birch Lumber
$ws.Cells.Item(89, 8).Value = 6271.143
$ws.Cells.Item(89, 9).Value = 4974.75
$ws.Cells.Item(89, 11).Value = 24873.75
$ws.Cells.Item(89, 13).Value = -19257.75

$ws = $wb.Worksheets.Item(5) # CUL
# Row 2: Pork Is a Salty Food / Table Salt
$ws.Cells.Item(2, 8).Value = 211.5
$ws.Cells.Item(2, 9).Value = 64.111115
$ws.Cells.Item(2, 11).Value = 384.66669
$ws.Cells.Item(2, 13).Value = -271.66669

# Row 21: Shy Is the Oyster / Raw Oyster
$ws.Cells.Item(21, 8).Value = 390
$ws.Cells.Item(21, 9).Value = 390
$ws.Cells.Item(21, 11).Value = 1170
$ws.Cells.Item(21, 13).Value = -997

# Row 22: A Total Nut Job / Walnut Bread
$ws.Cells.Item(22, 8).Value = 100
$ws.Cells.Item(22, 9).Value = 100
$ws.Cells.Item(22, 11).Value = 300
$ws.Cells.Item(22, 13).Value = -131

# Row 23: Sweet Smell of Success / Lavender Oil
$ws.Cells.Item(23, 8).Value = 474.22223
$ws.Cells.Item(23, 9).Value = 468
$ws.Cells.Item(23, 11).Value = 1404
$ws.Cells.Item(23, 13).Value = -1169

# Row 25: Flakes for Friends / Apple Tart
$ws.Cells.Item(25, 8).Value = 1007.25
$ws.Cells.Item(25, 9).Value = 469.33334
$ws.Cells.Item(25, 10).Value = 2621
$ws.Cells.Item(25, 11).Value = 1408.00002
$ws.Cells.Item(25, 12).Value = 7863
$ws.Cells.Item(25, 13).Value = -1239.00002
$ws.Cells.Item(25, 14).Value = -8201

# Row 27: Brain Food / Walnut Bread
$ws.Cells.Item(27, 8).Value = 100
$ws.Cells.Item(27, 9).Value = 100
$ws.Cells.Item(27, 11).Value = 300
$ws.Cells.Item(27, 13).Value = -198

# Row 30: Picnic Panic / Apple Tart
$ws.Cells.Item(30, 8).Value = 1007.25
$ws.Cells.Item(30, 9).Value = 469.33334
$ws.Cells.Item(30, 10).Value = 2621
$ws.Cells.Item(30, 11).Value = 1408.00002
$ws.Cells.Item(30, 12).Value = 7863
$ws.Cells.Item(30, 13).Value = -1306.00002
$ws.Cells.Item(30, 14).Value = -8067

$ws = $wb.Worksheets.Item(6) # GSM
# Row 101: Best-laid Planispheres / Dual-plated Durium Planisphere
$ws.Cells.Item(101, 8).Value = 49650
$ws.Cells.Item(101, 10).Value = 49650
$ws.Cells.Item(101, 12).Value = 49650
$ws.Cells.Item(101, 14).Value = -56140

# Row 113: Copious Crystal Cannons / Manasilver Nugget
$ws.Cells.Item(113, 8).Value = 1684.6
$ws.Cells.Item(113, 10).Value = 1853
$ws.Cells.Item(113, 12).Value = 1853
$ws.Cells.Item(113, 14).Value = -6193

$ws = $wb.Worksheets.Item(7) # LTW
# Row 16: Saddle Sore / Hard Leather
$ws.Cells.Item(16, 8).Value = 715.9167
$ws.Cells.Item(16, 9).Value = 469.2
$ws.Cells.Item(16, 10).Value = 1949.5
$ws.Cells.Item(16, 11).Value = 469.2
$ws.Cells.Item(16, 12).Value = 1949.5
$ws.Cells.Item(16, 13).Value = -299.2
$ws.Cells.Item(16, 14).Value = -2289.5

# Row 122: Hell on Leather / Gaja Leather
$ws.Cells.Item(122, 8).Value = 3945.4546
$ws.Cells.Item(122, 9).Value = 4050.125
$ws.Cells.Item(122, 10).Value = 3666.3333
$ws.Cells.Item(122, 11).Value = 12150.375
$ws.Cells.Item(122, 12).Value = 10998.9999
$ws.Cells.Item(122, 13).Value = -9700.375
$ws.Cells.Item(122, 14).Value = -15898.9999

$ws = $wb.Worksheets.Item(8) # WVR
# Row 81: Where the Dragonflies, the Net Catches / Crawler Silk
$ws.Cells.Item(81, 8).Value = 2212.8
$ws.Cells.Item(81, 9).Value = 1317.4546
$ws.Cells.Item(81, 11).Value = 2634.9092
$ws.Cells.Item(81, 13).Value = -1573.9092

# Row 84: To Kill a Dragon on Nameday (L) / Crawler Silk
$ws.Cells.Item(84, 8).Value = 2212.8
$ws.Cells.Item(84, 9).Value = 1317.4546
$ws.Cells.Item(84, 11).Value = 13174.546
$ws.Cells.Item(84, 13).Value = -7870.546

# Row 126: A Polished Purchase / Snow Linen
$ws.Cells.Item(126, 8).Value = 2098.5
$ws.Cells.Item(126, 9).Value = 2004.6875
$ws.Cells.Item(126, 11).Value = 6014.0625
$ws.Cells.Item(126, 13).Value = -3544.0625

# Row 133: Begin with the Basics / Snow Cotton Jacket
$ws.Cells.Item(133, 8).Value = 56000
$ws.Cells.Item(133, 10).Value = 60000
$ws.Cells.Item(133, 12).Value = 60000
$ws.Cells.Item(133, 14).Value = -70120

# Row 135: In Line with Linen / Mountain Linen Cloak of Casting
$ws.Cells.Item(135, 8).Value = 83499.75
$ws.Cells.Item(135, 10).Value = 93000
$ws.Cells.Item(135, 12).Value = 93000
$ws.Cells.Item(135, 14).Value = -103140

# Row 141: Silk for Sunperch / Thunderyards Silk Coat of Casting
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 13).ClearContents()
